# Applies scheduled market-price / profit refresh to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Source data is non-formula (externally fetched current market prices + derived profit columns),
# so each touched cell is written as a literal value, matching the upstream diff exactly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2833.1667
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 2999.75
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2999.75
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -3349.75
$ws.Range("H123").Value = 39110
$ws.Range("J123").Value = 39110
$ws.Range("L123").Value = 39110
$ws.Range("N123").Value = -48910
$ws.Range("H137").Value = 1870.1666
$ws.Range("I137").Value = 1525.5
$ws.Range("K137").Value = 4576.5
$ws.Range("M137").Value = -2026.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3192.6
$ws.Range("I32").Value = 2029.7759
$ws.Range("K32").Value = 2029.7759
$ws.Range("M32").Value = -1742.7759
$ws.Range("H61").Value = 1941.3684
$ws.Range("I61").Value = 1294.9117
$ws.Range("J61").Value = 7436.25
$ws.Range("K61").Value = 1294.9117
$ws.Range("L61").Value = 7436.25
$ws.Range("M61").Value = -1082.9117
$ws.Range("N61").Value = -7860.25
$ws.Range("H132").Value = 1494.7
$ws.Range("I132").Value = 1243.9615
$ws.Range("J132").Value = 3124.5
$ws.Range("K132").Value = 3731.8845
$ws.Range("L132").Value = 9373.5
$ws.Range("M132").Value = -1201.8845
$ws.Range("N132").Value = -14433.5
$ws.Range("H136").Value = 1941.3684
$ws.Range("I136").Value = 1294.9117
$ws.Range("J136").Value = 7436.25
$ws.Range("K136").Value = 3884.7351
$ws.Range("L136").Value = 22308.75
$ws.Range("M136").Value = -1334.7351
$ws.Range("N136").Value = -27408.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 17296
$ws.Range("J58").Value = 17296
$ws.Range("L58").Value = 17296
$ws.Range("N58").Value = -17884
$ws.Range("H61").Value = 17000
$ws.Range("J61").Value = 17000
$ws.Range("L61").Value = 17000
$ws.Range("N61").Value = -17626
$ws.Range("H86").Value = 1616.5238
$ws.Range("I86").Value = 1689.9
$ws.Range("J86").Value = 1549.8182
$ws.Range("K86").Value = 1689.9
$ws.Range("L86").Value = 1549.8182
$ws.Range("M86").Value = -566.9000000000001
$ws.Range("N86").Value = -3795.8182
$ws.Range("H89").Value = 1616.5238
$ws.Range("I89").Value = 1689.9
$ws.Range("J89").Value = 1549.8182
$ws.Range("K89").Value = 8449.5
$ws.Range("L89").Value = 7749.090999999999
$ws.Range("M89").Value = -2833.5
$ws.Range("N89").Value = -18981.091
$ws.Range("H107").Value = 2338.2144
$ws.Range("I107").Value = 1974.9
$ws.Range("K107").Value = 1974.9
$ws.Range("M107").Value = -54.90000000000009

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999.75
$ws.Range("I16").Value = 959.8
$ws.Range("K16").Value = 959.8
$ws.Range("M16").Value = -672.8
$ws.Range("H58").Value = 1228.6786
$ws.Range("I58").Value = 1049.6428
$ws.Range("J58").Value = 1407.7142
$ws.Range("K58").Value = 1049.6428
$ws.Range("L58").Value = 1407.7142
$ws.Range("M58").Value = -846.6428000000001
$ws.Range("N58").Value = -1813.7142
$ws.Range("H105").Value = 693.2
$ws.Range("I105").Value = 791.5
$ws.Range("K105").Value = 791.5
$ws.Range("M105").Value = 955.5
$ws.Range("H113").Value = 999.75
$ws.Range("I113").Value = 959.8
$ws.Range("K113").Value = 959.8
$ws.Range("M113").Value = 1210.2
$ws.Range("H132").Value = 1502.4103
$ws.Range("I132").Value = 997.40625
$ws.Range("J132").Value = 3811
$ws.Range("K132").Value = 2992.21875
$ws.Range("L132").Value = 11433
$ws.Range("M132").Value = -462.21875
$ws.Range("N132").Value = -16493
$ws.Range("H134").Value = 692.4545000000001
$ws.Range("I134").Value = 696.0526
$ws.Range("J134").Value = 669.6667
$ws.Range("K134").Value = 2088.1578
$ws.Range("L134").Value = 2009.0001
$ws.Range("M134").Value = 446.8422
$ws.Range("N134").Value = -7079.0001
$ws.Range("H136").Value = 1228.6786
$ws.Range("I136").Value = 1049.6428
$ws.Range("J136").Value = 1407.7142
$ws.Range("K136").Value = 3148.9284
$ws.Range("L136").Value = 4223.142599999999
$ws.Range("M136").Value = -598.9284000000002
$ws.Range("N136").Value = -9323.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3099.8333
$ws.Range("I3").Value = 1350
$ws.Range("K3").Value = 4050
$ws.Range("M3").Value = -3938
$ws.Range("H33").Value = 79.3
$ws.Range("I33").Value = 121.2
$ws.Range("J33").Value = 37.4
$ws.Range("K33").Value = 727.2
$ws.Range("L33").Value = 224.4
$ws.Range("M33").Value = -444.2
$ws.Range("N33").Value = -790.4
$ws.Range("H131").Value = 782.0599999999999
$ws.Range("I131").Value = 596.8
$ws.Range("J131").Value = 791.81055
$ws.Range("K131").Value = 1790.4
$ws.Range("L131").Value = 2375.43165
$ws.Range("M131").Value = 3249.6
$ws.Range("N131").Value = -12455.43165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5556.857
$ws.Range("I70").Value = 5979.6
$ws.Range("K70").Value = 5979.6
$ws.Range("M70").Value = -5709.6
$ws.Range("H73").Value = 5556.857
$ws.Range("I73").Value = 5979.6
$ws.Range("K73").Value = 5979.6
$ws.Range("M73").Value = -5043.6
$ws.Range("H102").Value = 2652.7856
$ws.Range("I102").Value = 3149.2856
$ws.Range("K102").Value = 3149.2856
$ws.Range("M102").Value = -1527.2856
$ws.Range("H126").Value = 50662.047
$ws.Range("I126").Value = 3427.4666
$ws.Range("K126").Value = 10282.3998
$ws.Range("M126").Value = -7812.399800000001
$ws.Range("H132").Value = 3046.628
$ws.Range("I132").Value = 2622.9412
$ws.Range("J132").Value = 4647.222
$ws.Range("K132").Value = 7868.823600000001
$ws.Range("L132").Value = 13941.666
$ws.Range("M132").Value = -5338.823600000001
$ws.Range("N132").Value = -19001.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1380.1
$ws.Range("J46").Value = 1400.125
$ws.Range("L46").Value = 1400.125
$ws.Range("N46").Value = -1776.125
$ws.Range("H55").Value = 185.88889
$ws.Range("I55").Value = 145
$ws.Range("J55").Value = 218.6
$ws.Range("K55").Value = 145
$ws.Range("L55").Value = 216.18182
$ws.Range("M55").Value = 28
$ws.Range("N55").Value = -564.6
$ws.Range("H132").Value = 1839.8148
$ws.Range("I132").Value = 1349.3846
$ws.Range("K132").Value = 4048.1538
$ws.Range("M132").Value = -1518.1538
$ws.Range("H136").Value = 5032.467
$ws.Range("I136").Value = 4312.5713
$ws.Range("K136").Value = 12937.7139
$ws.Range("M136").Value = -10387.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 74500
$ws.Range("J80").Value = 74500
$ws.Range("L80").Value = 74500
$ws.Range("N80").Value = -76496
$ws.Range("H83").Value = 74500
$ws.Range("J83").Value = 74500
$ws.Range("L83").Value = 223500
$ws.Range("N83").Value = -233484
$ws.Range("H132").Value = 2388.8518
$ws.Range("I132").Value = 2075.2856
$ws.Range("K132").Value = 6225.8568
$ws.Range("M132").Value = -3695.8568
$ws.Range("H136").Value = 3340.1072
$ws.Range("J136").Value = 3498.6
$ws.Range("L136").Value = 10495.8
$ws.Range("N136").Value = -15595.8

Write-Output "Updated 182 cells across 8 sheets"
